# PAS-6576: fixed VIN upload files for choice and select products,
# swap BI/PD/UM/MP symbol values (C <-> A) between row 2 and row 4,
# and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BI_SYMBOL/PD_SYMBOL/UM_SYMBOL/MP_SYMBOL change from "A" to "C"
$ws.Range("AE2").Value = "C"
$ws.Range("AF2").Value = "C"
$ws.Range("AG2").Value = "C"
$ws.Range("AH2").Value = "C"

# Row 4: BI_SYMBOL/PD_SYMBOL/UM_SYMBOL/MP_SYMBOL change from "C" to "A"
$ws.Range("AE4").Value = "A"
$ws.Range("AF4").Value = "A"
$ws.Range("AG4").Value = "A"
$ws.Range("AH4").Value = "A"

# Update the saved selection/active cell shown in the sheet view
$ws.Range("I13").Select()
